$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 69005
$ws.Range("E2").Value = 870
$ws.Range("F2").Value = 870
$ws.Range("G2").Value = 490
$ws.Range("H2").Value = 409
$ws.Range("I2").Value = 410
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 36582
$ws.Range("L2").Value = 24114
$ws.Range("M2").Value = 12468
$ws.Range("N2").Value = 11410
$ws.Range("O2").Value = 1057
$ws.Range("P2").Value = 343
$ws.Range("Q2").Value = 422
$ws.Range("R2").Value = -1742
$ws.Range("S2").Value = 1673
$ws.Range("T2").Value = 812
$ws.Range("U2").Value = -390
$ws.Range("V2").Value = 16848
$ws.Range("W2").Value = 1.26
$ws.Range("X2").Value = 0.59
$ws.Range("Y2").Value = 3.65
$ws.Range("Z2").Value = 1.12
$ws.Range("AA2").Value = 193.41
$ws.Range("AB2").Value = 3321.15
$ws.Range("AC2").Value = 5981
$ws.Range("AD2").Value = 11.55
$ws.Range("AE2").Value = 197350
$ws.Range("AF2").Value = 0.35
$ws.Range("AG2").Value = 2000
$ws.Range("AH2").Value = 2.89
$ws.Range("AI2").Value = 28.18
$ws.Range("AJ2").Value = 6860000

# Row 3
$ws.Range("D3").Value = 46143
$ws.Range("E3").Value = 317
$ws.Range("F3").Value = 357
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = -99
$ws.Range("K3").Value = 35098
$ws.Range("L3").Value = 22738
$ws.Range("M3").Value = 12360
$ws.Range("N3").Value = 11504
$ws.Range("O3").Value = 856
$ws.Range("P3").Value = 343
$ws.Range("Q3").Value = 1843
$ws.Range("R3").Value = -551
$ws.Range("S3").Value = -614
$ws.Range("T3").Value = 411
$ws.Range("U3").Value = 1432
$ws.Range("V3").Value = 16497
$ws.Range("W3").Value = 0.6899999999999999
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 0.87
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 183.96
$ws.Range("AB3").Value = 3314.34
$ws.Range("AC3").Value = 1454
$ws.Range("AD3").Value = 42.22
$ws.Range("AE3").Value = 198973
$ws.Range("AF3").Value = 0.31
$ws.Range("AG3").Value = 2600
$ws.Range("AH3").Value = 4.23
$ws.Range("AI3").Value = 150.7
$ws.Range("AJ3").Value = 6860000

# Row 4
$ws.Range("D4").Value = 39959
$ws.Range("E4").Value = 111
$ws.Range("F4").Value = 111
$ws.Range("G4").Value = -338
$ws.Range("H4").Value = -379
$ws.Range("I4").Value = -265
$ws.Range("J4").Value = -114
$ws.Range("K4").Value = 35442
$ws.Range("L4").Value = 23610
$ws.Range("M4").Value = 11832
$ws.Range("N4").Value = 11089
$ws.Range("O4").Value = 743
$ws.Range("P4").Value = 343
$ws.Range("Q4").Value = 282
$ws.Range("R4").Value = -105
$ws.Range("S4").Value = -167
$ws.Range("T4").Value = 457
$ws.Range("U4").Value = -175
$ws.Range("V4").Value = 16632
$ws.Range("W4").Value = 0.28
$ws.Range("X4").Value = -0.95
$ws.Range("Y4").Value = -2.35
$ws.Range("Z4").Value = -1.08
$ws.Range("AA4").Value = 199.53
$ws.Range("AB4").Value = 3193.59
$ws.Range("AC4").Value = -3869
$ws.Range("AD4").Value = -16.21
$ws.Range("AE4").Value = 191790
$ws.Range("AF4").Value = 0.33
$ws.Range("AG4").Value = 2000
$ws.Range("AH4").Value = 3.19
$ws.Range("AI4").Value = -43.57
$ws.Range("AJ4").Value = 6860000

# Row 5
$ws.Range("D5").Value = 44082
$ws.Range("E5").Value = 937
$ws.Range("F5").Value = 937
$ws.Range("G5").Value = 1230
$ws.Range("H5").Value = 868
$ws.Range("I5").Value = 844
$ws.Range("J5").Value = 23
$ws.Range("K5").Value = 34601
$ws.Range("L5").Value = 22091
$ws.Range("M5").Value = 12510
$ws.Range("N5").Value = 11743
$ws.Range("O5").Value = 767
$ws.Range("P5").Value = 343
$ws.Range("Q5").Value = -1290
$ws.Range("R5").Value = 1942
$ws.Range("S5").Value = -460
$ws.Range("T5").Value = 579
$ws.Range("U5").Value = -1868
$ws.Range("V5").Value = 16127
$ws.Range("W5").Value = 2.13
$ws.Range("X5").Value = 1.97
$ws.Range("Y5").Value = 7.39
$ws.Range("Z5").Value = 2.48
$ws.Range("AA5").Value = 176.58
$ws.Range("AB5").Value = 3403.22
$ws.Range("AC5").Value = 12306
$ws.Range("AD5").Value = 4.62
$ws.Range("AE5").Value = 171187
$ws.Range("AF5").Value = 0.33
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 3.51
$ws.Range("AI5").Value = 13.7
$ws.Range("AJ5").Value = 6860000

# Row 6
$ws.Range("D6").Value = 46302
$ws.Range("E6").Value = 140
$ws.Range("F6").Value = 140
$ws.Range("G6").Value = 1053
$ws.Range("H6").Value = 792
$ws.Range("I6").Value = 781
$ws.Range("K6").Value = 35654
$ws.Range("L6").Value = 22495
$ws.Range("M6").Value = 13159
$ws.Range("N6").Value = 12381
$ws.Range("P6").Value = 343
$ws.Range("Q6").Value = 1847
$ws.Range("R6").Value = -268
$ws.Range("S6").Value = -2215
$ws.Range("T6").Value = 521
$ws.Range("U6").Value = 1327
$ws.Range("V6").Value = 14030
$ws.Range("W6").Value = 0.3
$ws.Range("X6").Value = 1.71
$ws.Range("Y6").Value = 6.47
$ws.Range("Z6").Value = 2.26
$ws.Range("AA6").Value = 170.95
$ws.Range("AB6").Value = 3582.96
$ws.Range("AC6").Value = 11384
$ws.Range("AD6").Value = 4.9
$ws.Range("AE6").Value = 180483
$ws.Range("AF6").Value = 0.31
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 3.58
$ws.Range("AI6").Value = 14.81
$ws.Range("AJ6").Value = 6860000

# Clear rows 7-9 data columns (D through AI), keep A, B, C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
